$d = $word.ActiveDocument

# Locate the paragraph right after "LOB1018: Física I (Requisito fraco)"
# (an empty paragraph) and the paragraph holding the "© 2020 ... Creative
# Commons Attribution" footer text; together with the "Ver no Jupiter..."
# paragraph in between, these three paragraphs are removed so that the
# requisitos line flows straight into the (remaining) blank paragraph that
# precedes the page-break paragraph.
$count = $d.Paragraphs.Count
$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "*LOB1018*") {
        $startIndex = $i + 1
    }
    if ($text -like "*Creative Commons Attribution*") {
        $endIndex = $i
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    $rangeStart = $d.Paragraphs.Item($startIndex).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIndex).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
